# LED Control test 1
# Fill in row 10 of the Translation sheet with a new "LED" text entry,
# reusing the reserved "SingleUseId9" text id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Reset to the default/normal style so the new cells don't pick up a
# stray explicit style reference (matches the style-less cells already
# used by the other data rows in this table).
$ws.Range("B10:F10").Style = "Normal"

$ws.Range("B10").Value = "SingleUseId9"
$ws.Range("C10").Value = "ButtonDown"
$ws.Range("D10").Value = "Center"
$ws.Range("E10").Value = "LED"
$ws.Range("F10").Value = "LTR"
